# New crime data collected - weekly update for week ending 3/2/2025
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: copy both the format (style) and the value from a source cell onto
# a destination cell, using two PasteSpecial passes so the destination picks
# up exactly the same style index as the source (rather than Excel minting a
# brand-new style record), and the value/text-vs-number typing of the source.
function Copy-CellFull {
    param([string]$srcAddr, [string]$dstAddr)
    $src = $ws.Range($srcAddr)
    $dst = $ws.Range($dstAddr)
    $src.Copy() | Out-Null
    $dst.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    $src.Copy() | Out-Null
    $dst.PasteSpecial(-4163) | Out-Null  # xlPasteValues
    $excel.CutCopyMode = 0
}

# ---------------------------------------------------------------------
# Header text: volume number and report week dates
# ---------------------------------------------------------------------
$ws.Range("M8").Characters(23, 1).Text = "9"
$ws.Range("C9").Characters(27, 9).Text = "2/24/2025"
$ws.Range("C9").Characters(49, 6).Text = "3/2/2025"

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
Copy-CellFull "C18" "C15"
Copy-CellFull "C18" "D15"
Copy-CellFull "E22" "E15"
$ws.Range("N15").Value = -85

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 17
$ws.Range("H16").Value = 21.428571428571
$ws.Range("I16").Value = 36
$ws.Range("J16").Value = 32
$ws.Range("K16").Value = 12.5
$ws.Range("L16").Value = 33.333333333333
$ws.Range("M16").Value = -33.333333333333
$ws.Range("N16").Value = -90.163934426229

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 42.857142857142
$ws.Range("F17").Value = 36
$ws.Range("G17").Value = 26
$ws.Range("H17").Value = 38.461538461538
$ws.Range("I17").Value = 68
$ws.Range("J17").Value = 52
$ws.Range("K17").Value = 30.769230769230
$ws.Range("L17").Value = 11.475409836065
$ws.Range("M17").Value = 4.615384615384
$ws.Range("N17").Value = -47.692307692307

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
Copy-CellFull "C17" "C18"
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -80
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -40
$ws.Range("I18").Value = 12
$ws.Range("J18").Value = 19
$ws.Range("K18").Value = -36.842105263157
$ws.Range("L18").Value = -55.555555555555
$ws.Range("M18").Value = -67.567567567567
$ws.Range("N18").Value = -97.887323943662

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 50
$ws.Range("F19").Value = 35
$ws.Range("G19").Value = 39
$ws.Range("H19").Value = -10.256410256410
$ws.Range("I19").Value = 67
$ws.Range("J19").Value = 75
$ws.Range("K19").Value = -10.666666666666
$ws.Range("L19").Value = -40.178571428571
$ws.Range("M19").Value = -18.292682926829
$ws.Range("N19").Value = -59.638554216867

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 1
Copy-CellFull "C18" "D20"
Copy-CellFull "E22" "E20"
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = -12.5
$ws.Range("I20").Value = 13
$ws.Range("K20").Value = -35
$ws.Range("L20").Value = -27.777777777777
$ws.Range("M20").Value = -58.064516129032
$ws.Range("N20").Value = -96.813725490196

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 29
$ws.Range("E21").Value = 11.538461538461
$ws.Range("F21").Value = 103
$ws.Range("G21").Value = 99
$ws.Range("H21").Value = 4.040404040404
$ws.Range("I21").Value = 199
$ws.Range("J21").Value = 200
$ws.Range("K21").Value = -0.5
$ws.Range("L21").Value = -19.758064516129
$ws.Range("M21").Value = -26.838235294117
$ws.Range("N21").Value = -88.062387522495

# ---------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------
$ws.Range("L22").Value = -75
$ws.Range("M22").Value = -80

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 44
$ws.Range("D24").Value = 46
$ws.Range("E24").Value = -4.347826086956
$ws.Range("F24").Value = 115
$ws.Range("G24").Value = 192
$ws.Range("H24").Value = -40.104166666666
$ws.Range("I24").Value = 237
$ws.Range("J24").Value = 357
$ws.Range("K24").Value = -33.613445378151
$ws.Range("L24").Value = -10.566037735849
$ws.Range("M24").Value = 25.396825396825

# ---------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 17
$ws.Range("D25").Value = 22
$ws.Range("E25").Value = -22.727272727272
$ws.Range("F25").Value = 46
$ws.Range("G25").Value = 112
$ws.Range("H25").Value = -58.928571428571
$ws.Range("I25").Value = 128
$ws.Range("J25").Value = 213
$ws.Range("K25").Value = -39.906103286385
$ws.Range("L25").Value = -11.111111111111

# ---------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = 27.272727272727
$ws.Range("F26").Value = 42
$ws.Range("G26").Value = 51
$ws.Range("H26").Value = -17.647058823529
$ws.Range("I26").Value = 86
$ws.Range("J26").Value = 101
$ws.Range("K26").Value = -14.851485148514
$ws.Range("L26").Value = 2.380952380952
$ws.Range("M26").Value = -28.333333333333

# ---------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------
Copy-CellFull "C18" "C27"
Copy-CellFull "C18" "D27"
Copy-CellFull "E22" "E27"
$ws.Range("L27").Value = -50

# ---------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------
Copy-CellFull "C18" "C28"
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = -87.5
$ws.Range("J28").Value = 13
$ws.Range("K28").Value = -53.846153846153

# ---------------------------------------------------------------------
# Row 31 - Hate Crimes
# ---------------------------------------------------------------------
Copy-CellFull "C18" "D31"
Copy-CellFull "E22" "E31"

Write-Output "done"
